$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename teacher speaker tags from "T/R1" to "T1" (covers both "T/R1" and "T/R1:")
$ws.Cells.Replace("T/R1", "T1")

# Rename "Students" to "SS" (covers "Students", "Students:", and
# "3 - Getting Students to Relate" -> "3 - Getting SS to Relate")
$ws.Cells.Replace("Students", "SS")
